# Review_223.docx edit: update title date, swap in the new paper review
# text (LOOKAHEAD DECODING instead of CLLMs), and drop the trailing
# paragraphs that no longer apply, replacing the final link.
#
# Find/Replace (rather than a direct Range.Text assignment) is used so
# that Word recomputes the run's xml:space="preserve" flag from the new
# replacement text instead of inheriting it from the old run.

$d = $word.ActiveDocument

# --- 1) Title line: date 14.06.24 -> 13.06.24 ---------------------------
$d.Content.Find.Execute(
    "⚡️🚀המאמר היומי של מייק 14.06.24:⚡️🚀", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "⚡️🚀המאמר היומי של מייק 13.06.24:⚡️🚀", 2)

# --- 2) Paper title -------------------------------------------------------
$d.Content.Find.Execute(
    "CLLMs: Consistency Large Language Models", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Break the Sequential Dependency of LLM Inference Using LOOKAHEAD DECODING", 2)

# --- 3) Body paragraph 3 ---------------------------------------------------
$d.Content.Find.Execute(
    "בשתי הסקירות הקודמות(כדאי שתעברו עליהם כי נתתי שם קצת הסברים) דיברנו על שיטות איטרטיביות מקבילות לדגימה ממודלי שפה. השיטות האלו מבוססות על שיטות יאקובי או (Gauss-Seidel (GS. השיטות האלו מתחילות מכמות מסוימת n של טוקנים שנדגמים באקראי (או בצורה קצת יותר מושכלת) ואז מעדכנים טוקנים אלו בבת אחת באיטרציות עד שתנאי עצירה מתקיים(התכנסות). תנאי העצירה כאן הוא בד״כ שוויון בין הפלטים של איטרציות עוקבות. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "זוכרים את המאמר שסקרנו קצרות אתמול שהציע גישה איטרטיבית לפתרון מקבילי של מערכות משוואות לא לינאריות. אחת הדוגמאות של פתרון מערכות משוואות כאלו היא גנרוט טקסט ממודלי שפה כאשר כל טוקן נבחר בתור argmax של התפלגות הטוקן בהינתן הטוקנים הקודמים (המופק באמצעות השכבה האחרונה של מודל השפה).", 2)

# --- 4) Body paragraph 4 ---------------------------------------------------
$d.Content.Find.Execute(
    "מובן שאנו מעוניינים לסיים את התהליך במשמעות פחות איטרציות ממספר הטוקנים שאנו חוזים בו זמנית (ד״א ניתן להראות נדרשות לכל היותר ח איטרציות עד ההתכנסות). ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "יש בגדול שתי שיטות איטרטיביות שניתן לרתום אותן לדגימה יעילה יותר ממודלי שפה: יעקובי וגאוס-סיידל. שתי השיטות מתחילות מניחוש אקראי של כמה טוקנים בהינתן ההקשר ואז מאפטמים אותם על פתרון איטרטיבי של מערכת המשוואות עם argmax (ששקול לגנרוט). אפשר די בקלות לראות שבגלל שהמשוואות הן אוטורגרסיביות שיטות אלו לא יכולות להתכנס ביותר מ n איטרציות (מספר הטוקנים הנדגמים עם שיטה) ולפעמים אפשר להספיק פחות (נציין כי כל איטרציה דורשת קצת יותר משאבי החישוב).", 2)

# --- 5) Body paragraph 5 (new text keeps trailing space) -------------------
$d.Content.Find.Execute(
    "שימו לב שמהלך האימון של מודלי שפה מותאם לשיטת הדגימה האוטו-רגרסיביות כאשר בוחרים טוקן בעל הסתסברות הגבוה ביותר ביהנתן הטוקנים הקודמים. אולם עכשיו אנו דוגמים בצורה אחרת ואולי ניתן להתחשב בזה במהלך האימון. כלומר במהלך האימון אשכרה דוגמים עם השיטה הזו (השילוב של יאקובי ו- GS).",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "הבעיה עם השימוש הנאיבי בשיטה הוא שהרווח הממוצע על פני דגימה אוטורגרסיבית סטנדרטית ממודלי שפה הוא לא גדול ועומד על פחות מ 1.1 האצת קצב גנרוט. ", 2)

# --- 6) Body paragraph 6 ---------------------------------------------------
$d.Content.Find.Execute(
    "וזה בדיוק מה שנסקור אותו היום עושה. המחברים מוסיפים עוד איבר ללוס הרגיל של מודלי שפה (הממקסם את הנראות המירבית של הדאטה). מטרת האיבר הזה היא לגרום למזעור של מספר האיטרציות עד להתכנסות של הדגימה האיטרטיבית. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "המאמר מציע שכלול לשיטה הנאיבית ומציע לשמור בזכרון את הטוקנים של כמה איטרציות האחרונות. במקרה אם והיא מוצאת בזכרון זה תת-סדרת טוקנים שבה הטוקן הראשון זהה לטוקן הראשון ״הנכון״ של האיטרציה(באיטרציה i טוקן i וקודמיו נחזים נכון) אנו לוקחים תת סדרה זו ומציבים אותו במקום מה שנחזה באיטרציה האחרונה.", 2)

# --- 7) Body paragraph 7 ---------------------------------------------------
$d.Content.Find.Execute(
    "המחברים בחנו שתי אופציות לאיבר הזה:", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "זה מאפשר להקטין את כמות האיטרציות עוד טיפה", 2)

# --- 8) Final link paragraph: replace text, then delete the old trailing --
#        paragraphs that are no longer part of the review.
$d.Content.Find.Execute(
    "מזעור של מרחק (KL הפוך לדעתי אך לא צללתי לעומק) בין התפלגויות הטוקנים בנקודת ההתכנסות לבין התפלגויות טוקנים במהלך הדגימה האיטרטיבית (דוגמים האיטרציות באקראי).",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "https://arxiv.org/pdf/2402.02057", 2)

$lastIndex = $d.Paragraphs.Count
$startPara = $d.Paragraphs.Item(9)
$endPara = $d.Paragraphs.Item($lastIndex)
$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
